# Remove column M ("alcohol data" measurement column) from Sheet1.
# The old column N shifts left to become the new column M, matching the
# commit "remove column from alcohol data".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns.Item(13).Delete() | Out-Null

# Reflect the new active selection left behind after the column delete.
$ws.Range("M1").Select() | Out-Null
